# Rename the two logo pictures that are duplicated across the document's
# "primary" and "first page" headers/footers.
#
#   headers (BTec_Logo-Orange):        image2.jpg -> image1.jpg
#   footers (PearsonLogo):             image1.png -> image2.png
#
# InlineShapes live in four separate stories (primary header, first-page
# header, primary footer, first-page footer). Section.Headers.Item(1) /
# Footers.Item(1) address the "primary" story, Item(2) addresses the
# "first page" story.
#
# Re-fetching the InlineShape through its own .Range.InlineShapes.Item(1)
# (instead of reusing the originally-fetched reference) keeps the COM
# handle correctly addressed for every story, including footers, before
# the .Name write.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$sec.Headers.Item(1).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1).Name = "image1.jpg"
$sec.Headers.Item(2).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1).Name = "image1.jpg"

$sec.Footers.Item(1).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1).Name = "image2.png"
$sec.Footers.Item(2).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1).Name = "image2.png"
